$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Correct the marks: total correct answers and total max marks
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 20
$ws.Range("E12").Value = "20/140"
